$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "sequence/run_0711_5_0718_7/"

for ($r = 2; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -and $val.ToString().StartsWith($prefix)) {
        $cell.Value2 = $val.ToString().Substring($prefix.Length)
    }
}

$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 13
$win.ScrollColumn = 1

$ws.Range("F46").Select()
